# Update the division exercises contained in the document's single table.
# Every "NN÷N=" expression in the table is unique at the time it is looked
# up, EXCEPT for "63÷6=", which is both the original text of one cell and
# the new text that another cell is changed to. To avoid the second
# replacement accidentally matching the text freshly produced by the
# first, the pair is reordered so the cell that currently reads "63÷6="
# is changed away from it before another cell is changed into "63÷6=".

$d = $word.ActiveDocument

$edits = @(
    @("37÷3=", "44÷2="),
    @("53÷4=", "30÷2="),
    @("29÷2=", "59÷6="),
    @("63÷6=", "27÷8="),
    @("43÷7=", "63÷6="),
    @("98÷7=", "63÷8="),

    @("84÷2=", "14÷5="),
    @("65÷2=", "18÷5="),
    @("26÷9=", "93÷3="),
    @("41÷9=", "24÷4="),
    @("38÷2=", "13÷4="),

    @("12÷2=", "16÷9="),
    @("95÷8=", "35÷6="),
    @("62÷6=", "41÷8="),
    @("35÷4=", "30÷4="),

    @("65÷3=", "29÷5="),
    @("56÷2=", "44÷2="),
    @("50÷5=", "97÷4="),
    @("17÷8=", "35÷8="),
    @("95÷2=", "83÷3="),

    @("66÷8=", "54÷8="),
    @("18÷8=", "49÷8="),
    @("72÷5=", "44÷8="),
    @("70÷2=", "87÷3="),
    @("58÷7=", "73÷6=")
)

foreach ($e in $edits) {
    $old = $e[0]
    $new = $e[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2) | Out-Null
}
